$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.679.28'
$ws.Range('E2').Value = '  +0.04%  '

$ws.Range('D3').Value = '1.850.15'
$ws.Range('E3').Value = '  +0.50%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.67%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '313.10'
$ws.Range('E5').Value = '  -0.78%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.62%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4234'
$ws.Range('E7').Value = '  +0.36%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3645'
$ws.Range('E8').Value = '  +0.50%  '

$ws.Range('E9').Value = '  -0.14%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07293'
$ws.Range('E10').Value = '  +1.10%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.8774'
$ws.Range('E11').Value = '  -2.17%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '20.64'
$ws.Range('E12').Value = '  +0.37%  '

$ws.Range('D13').Value = '1.826.49'
$ws.Range('E13').Value = '  +0.23%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.327'

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.527'
$ws.Range('E15').Value = '  -0.28%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.06856'
$ws.Range('E16').Value = '  +0.45%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.000'
$ws.Range('E17').Value = '  -0.80%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '79.72'
$ws.Range('E18').Value = '  +2.46%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000008931'
$ws.Range('E19').Value = '  +0.22%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  -0.51%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '15.34'
$ws.Range('E21').Value = '  +0.12%  '

$ws.Range('D22').Value = '27.678.03'
$ws.Range('E22').Value = '  +0.11%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.983'
$ws.Range('E23').Value = '  +0.76%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '10.35'
$ws.Range('E24').Value = '  -4.16%  '

$ws.Range('D25').Value = '2.073.96'
$ws.Range('E25').Value = '  +2.26%  '

$ws.Range('E26').Value = '  -2.68%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '154.17'
$ws.Range('E27').Value = '  +0.54%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.83'
$ws.Range('E28').Value = '  +3.86%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '122.14'
$ws.Range('E29').Value = '  +10.11%  '

$ws.Range('E30').Value = '  -0.49%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.873'
$ws.Range('E31').Value = '  +15.61%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.08852'
$ws.Range('E32').Value = '  -0.32%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.7670'
$ws.Range('E33').Value = '  -0.86%  '

$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.542'
$ws.Range('E34').Value = '  +0.82%  '

$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.971'
$ws.Range('E35').Value = '  +1.08%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.106'
$ws.Range('E36').Value = '  +3.26%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.9997'
$ws.Range('E37').Value = '  -0.72%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.095'
$ws.Range('E38').Value = '  +1.24%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.05358'
$ws.Range('E39').Value = '  +0.60%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01933'
$ws.Range('E40').Value = '  +0.99%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.819'
$ws.Range('E41').Value = '  -4.56%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.896'
$ws.Range('E42').Value = '  +3.10%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.5084'
$ws.Range('E43').Value = '  +0.35%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.1648'
$ws.Range('E44').Value = '  +0.76%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.331'
$ws.Range('E45').Value = '  +1.59%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.06538'
$ws.Range('E46').Value = '  -1.48%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.34'
$ws.Range('E47').Value = '  +1.17%  '

$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '105.55'
$ws.Range('E48').Value = '  +0.33%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.4686'
$ws.Range('E49').Value = '  -0.18%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.9996'
$ws.Range('E50').Value = '  -0.73%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.625'
$ws.Range('E51').Value = '  +0.26%  '
